# Noise Removal in pipeline
# Update the Leaf_size(px) (col F) and Lateral_root_count (col G) values
# for rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2004
$ws.Range("G2").Value = 3

$ws.Range("F3").Value = 2174
$ws.Range("G3").Value = 4

$ws.Range("F4").Value = 2291

$ws.Range("F5").Value = 1972

$ws.Range("F6").Value = 1352
$ws.Range("G6").Value = 4
